# Update "想去人数" (wish-to-go count) values in column F across sheets
# 展览 (Exhibition), 演出 (Performance), 全部类型 (All Types)

$wb = $excel.ActiveWorkbook

$updates = @{
    "展览" = @{
        "F2" = 639
        "F3" = 704
        "F4" = 945
        "F5" = 720
        "F6" = 838
        "F7" = 402
        "F8" = 601
        "F9" = 131
        "F10" = 1208
        "F11" = 635
        "F13" = 508
        "F15" = 13
        "F16" = 484
        "F17" = 355
        "F18" = 59
        "F20" = 555
        "F21" = 77
        "F22" = 578
        "F24" = 758
    }
    "演出" = @{
        "F2" = 86
        "F9" = 222
        "F13" = 94
    }
    "全部类型" = @{
        "F3" = 86
        "F4" = 639
        "F7" = 704
        "F8" = 945
        "F9" = 720
        "F10" = 838
        "F11" = 402
        "F12" = 601
        "F13" = 131
        "F14" = 1208
        "F15" = 635
        "F19" = 508
        "F22" = 13
        "F23" = 484
        "F25" = 355
        "F26" = 59
        "F28" = 222
        "F30" = 555
        "F33" = 94
        "F34" = 94
        "F35" = 77
        "F36" = 578
        "F38" = 758
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $cellMap = $updates[$sheetName]
    foreach ($cellRef in $cellMap.Keys) {
        $ws.Range($cellRef).Value = $cellMap[$cellRef]
    }
}
